# Update workbook with residential water heating load data.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "Water" sheet
$ws2 = $wb.Worksheets.Item(2)   # "Air" sheet

# New normalized load profile values (A1:X1), identical on both sheets.
$vals = @(
  0.24682862601720568,
  0,
  0,
  0,
  0,
  0.060065141685888224,
  0,
  0.059144537942310671,
  0.51927141127240273,
  0.52382438939645537,
  0.13038529760314024,
  0.51398284812689254,
  0.016195395117100556,
  0,
  0.26090828123027648,
  0.089361092814150195,
  0.013040826386388669,
  0.2276205737539557,
  0.034640545575266177,
  0.20990934955450419,
  0.054831668975116259,
  0.38097232065841097,
  0.45580618727025773,
  0.061190209757303174
)

for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws1.Cells.Item(1, $i + 1).Value = $vals[$i]
    $ws2.Cells.Item(1, $i + 1).Value = $vals[$i]
}

# Update the selected cell on each sheet.
$null = $ws1.Range("C25").Select()
$null = $ws2.Range("C21").Select()

# Make the "Air" sheet the active / selected tab (was "Water").
$null = $ws2.Activate()
